# StoryCards.xlsx update
# 1. Two new StoryCards created (Doku & Vortrag)
# 2. Status/Akzeptanztest updates for several existing story cards

$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Row 14 (StoryCard ID 7): mark done + accepted ---
$ws.Range("B14").Value = "fertig"
$ws.Range("J14").Value = "akzeptiert"

# --- Row 20 (StoryCard ID 13): mark done + accepted ---
$ws.Range("B20").Value = "fertig"
$ws.Range("J20").Value = "akzeptiert"

# --- Row 21 (StoryCard ID 14): mark accepted ---
$ws.Range("J21").Value = "akzeptiert"

# --- Row 22 (StoryCard ID 15): mark accepted ---
$ws.Range("J22").Value = "akzeptiert"

# --- Row 23 (StoryCard ID 16): mark done + accepted, add actual effort + finish date ---
$ws.Range("B23").Value = "fertig"
$ws.Range("J23").Value = "akzeptiert"
$ws.Range("L23").Value = "2h 30min"
$ws.Range("M23").Value = 40830
$ws.Range("M22").Copy()
$ws.Range("M23").PasteSpecial($xlPasteFormats)

# --- New Row 24 (StoryCard ID 17): "Vortrag" ---
$ws.Range("A24").Value = 17
$ws.Range("B24").Value = "in Arbeit"
$ws.Range("C24").Value = "hoch"
$ws.Range("K24").Value = "4h"
$ws.Range("D24").Value = "Vortrag"
$ws.Range("E24").Value = 40840
$ws.Range("F24").Value = "Wiederschein"
$ws.Range("G24").Value = "alle"
$ws.Range("H24").Value = "Ausarbeitung"

# --- New Row 25 (StoryCard ID 18): "Dokumentation" ---
$ws.Range("A25").Value = 18
$ws.Range("B25").Value = "in Arbeit"
$ws.Range("C25").Value = "hoch"
$ws.Range("D25").Value = "Dokumentation"
$ws.Range("E25").Value = 40840
$ws.Range("F25").Value = "Wiederschein"
$ws.Range("G25").Value = "alle"
$ws.Range("H25").Value = "Ausarbeitung"
$ws.Range("K25").Value = "8h"

# Beschreibung text for the two new cards
$ws.Range("I24").Value = "Präsentation erstellen"
$ws.Range("I25").Value = "Dokumentation erstellen"

# Match date formatting used by the other rows in the Datum column (same format as E23)
$ws.Range("E23").Copy()
$ws.Range("E24:E25").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false
